# "Cellular for SMART and FLOW fixed."
#
# Main Info sheet: fix the Design / Converged-router / Backup_IP+mask /
# 4G+Cellular fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")
$ws.Activate()

# Design: BASE -> FLOW
$ws.Range("B7").Value = "FLOW"

# Converged router: TRUE -> FALSE
$ws.Range("B8").Value = $false

# Backup_IP+mask: dhcp -> 2.2.2.2/24
$ws.Range("B20").Value = "2.2.2.2/24"

# 4G+Cellular: FALSE -> TRUE
$ws.Range("B25").Value = $true

# Scroll the view down a row and move the active selection from D7 to D8
# (matches the saved view state in the workbook).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D8").Select()
